# "adding averages and more checks"
# - Update "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I) on Training Dashboard (LAST UPDATE moved
#   forward from 08-Sep-2025 to 16-Sep-2025, so the remaining-days counters in column H drop by 8).
# - Exam Dashboard: widen the COMMENTS column and change the per-row comment from "OK" to the more
#   descriptive "date is valid".
# - Re-colour the title banners and table headers on both sheets to bold white text (on their dark
#   fill), dropping the old oversized 14pt title font in favour of the normal 11pt size.

$wb = $excel.ActiveWorkbook

$training = $wb.Worksheets.Item("Training Dashboard")
$exam     = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------------
# 1. Training Dashboard: refresh "LAST UPDATE" date + recomputed "PERIOD TO EXPIRE"
# ---------------------------------------------------------------------------
$periodToExpire = @{
    3  = 605
    4  = 604
    5  = 604
    6  = 604
    7  = 604
    8  = 605
    9  = 608
    10 = 605
    11 = 608
    12 = 605
    13 = 608
    14 = -19618
    15 = 284
    16 = 309
    17 = 309
}

foreach ($row in 3..17) {
    $training.Range("H$row").Value = $periodToExpire[$row]
    # Keep this a literal text date (matching the sheet's existing style) rather than
    # letting Excel's autocoerce turn it into a real date serial number.
    $training.Range("I$row").Value = "'16-Sep-2025"
}

# ---------------------------------------------------------------------------
# 2. Exam Dashboard: wider COMMENTS column + clearer comment text
# ---------------------------------------------------------------------------
$exam.Columns.Item(5).ColumnWidth = 14.166666666666666

foreach ($row in 3..5) {
    $exam.Range("E$row").Value = "date is valid"
}

# ---------------------------------------------------------------------------
# 3. Header / title formatting: bold white text (drop the separate 14pt title font)
# ---------------------------------------------------------------------------
foreach ($ws in @($training, $exam)) {
    $titleRange = $ws.Range("A1")
    $titleRange.Font.Size = 11
    $titleRange.Font.Color = 16777215

    $lastCol = $ws.Cells.Item(2, $ws.UsedRange.Columns.Count).Column
    $headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $headerRange.Font.Size = 11
    $headerRange.Font.Color = 16777215
}
